$d = $word.ActiveDocument

# Locate the paragraph that introduces the block to remove: the one
# immediately following "LOQ4239: Administração e Organização I (Requisito
# fraco)". The three paragraphs to delete are:
#   1) an empty spacer paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "© 2020 . Contact: ..." footer paragraph
$count = $d.Paragraphs.Count
$startIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "LOQ4239*") {
        $startIndex = $i + 1
        break
    }
}

if ($startIndex -gt 0) {
    $endIndex = $startIndex + 2
    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($endIndex)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
